$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values in column F (dSF) for specific rows, per repull of data / mean calculation fix
$updates = @{
    5  = -1
    7  = -3
    8  = -4
    10 = -4
    12 = -4
    13 = -8
    17 = -2
    20 = -4
    23 = -3
    24 = -3
    27 = -4
    28 = -10
    31 = 0
    38 = 0
    39 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
